# Adding the changes we made on may 9th
$wb2 = $excel.ActiveWorkbook
$ws = $wb2.ActiveSheet

# Insert one new data row at row 2 (shifts existing rows 2-21 down to rows 3-22).
$ws.Rows("2:2").Insert()
# Excel's row Insert copies the formatting of the row above (the bold/bordered
# header) onto the new row; the new data row should carry no explicit style,
# matching the rest of the (unstyled) data rows.
$ws.Rows("2:2").ClearFormats()

# Populate the freshly-inserted row 2 with the new sample.
$row2 = New-Object 'object[,]' 1,3
$row2[0,0] = 0.7064247653999316
$row2[0,1] = 0.3879705256011861
$row2[0,2] = -1.276251717872426
$ws.Range("A2:C2").Value = $row2

# Append the new trailing rows (23-31) with additional samples.
$newRows = @(
    @(3.552849229822288, -5.358195619534956, -0.1442216574237578),
    @(4.79510967652812, 16.44486069921317, -4.148945381193625),
    @(-0.9187659382214614, -0.2956210344576297, 3.151463126168018),
    @(-5.69272972968636, 3.907055351334684, 4.745034019354134),
    @(-0.1373755151245123, 11.3405332758948, -3.450026544822681),
    @(-1.20526529810765, -4.696971940510186, 1.208450563062909),
    @(8.444541233142751, -2.498347297840359, -2.655560967885902),
    @(-0.723370986541445, 0.8442579066087497, -2.731481316125937),
    @(-3.363093618511545, 1.143067340257809, 1.45013582887988)
)

$rowCount = $newRows.Count
$arr = New-Object 'object[,]' $rowCount,3
for ($i = 0; $i -lt $rowCount; $i++) {
    $vals = $newRows[$i]
    for ($j = 0; $j -lt 3; $j++) {
        $arr[$i,$j] = $vals[$j]
    }
}

$startRow = 23
$endRow = $startRow + $rowCount - 1
$ws.Range("A" + $startRow + ":C" + $endRow).Value = $arr
